$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: Remove the "Meta description: ..." paragraph that used to
# sit right under the "Play Feng Fu Online Slot Game for Free" heading.
# ---------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$null = $metaPara.Range.Delete()

# ---------------------------------------------------------------------
# Change 2: Insert a new bold "Play Feng Fu Online Slot Game for Free"
# paragraph right before the final "Prompt: ..." paragraph.
# ---------------------------------------------------------------------
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($n)
$null = $lastPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs($n)
$ooxml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Feng Fu Online Slot Game for Free</w:t></w:r></w:p>'
$null = $newPara.Range.InsertXML($ooxml)

# ---------------------------------------------------------------------
# Change 3: Replace the text of the final "Prompt: ..." paragraph with
# the new meta-description text (keeping its italic formatting).
# ---------------------------------------------------------------------
$oldText = "Prompt: Create a cartoon-style feature image for the game " + [char]34 + "Feng Fu" + [char]34 + " featuring a happy Maya warrior with glasses. The image should include a background of traditional Chinese elements such as the Great Wall of China or the bamboo scrolls that are depicted in the game. The Maya warrior should be wearing a playful expression and should be holding a gong or plate, which is a special symbol in the game. The overall tone of the image should be colorful and vibrant, in line with the game's theme of luck and fortune."
$newText = "Explore the theme of luck and wealth in Feng Fu online slot game, using Chinese symbolism and ideograms. Play for free and enjoy high value wins."

$null = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false,
                                 $true, 1, $false, $newText, 2)
